$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values for rows 2-11: vlam (C) 0.05 -> 0.08, vlag (G) 4 -> 9, mlos (H) 12 -> 13, vlos (I) 50 -> 154
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = 0.08
    $ws.Cells.Item($r, 7).Value = 9
    $ws.Cells.Item($r, 8).Value = 13
    $ws.Cells.Item($r, 9).Value = 154
}

# Update default column width (sheetFormatPr defaultColWidth 12.0390625 -> 12.0546875)
$ws.StandardWidth = 12.0546875

# Update selection / active cell to I18
$ws.Range("I18").Select()
